# Applies the "Updated cryptos list" refresh (Sun Sep 24 22:37:06 UTC 2023).
# Column layout: B=Coin, C=Link, D=Price, E=Volume(1h). Column A (rank) is untouched.
#
# Price-column values that look like plain numbers (e.g. "1.00", "210.73") are written
# with a leading apostrophe so Excel keeps them as literal text (matching the workbook's
# inline-string cells, e.g. "1.00" must stay "1.00" and not collapse to 1). The apostrophe
# marker is stripped from the stored value automatically; ".Style = 'Normal'" immediately
# afterwards clears the resulting quote-prefix formatting so the cell style is unchanged.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: Bitcoin
$ws.Range("D2").Value = "26.623.56"
$ws.Range("E2").Value = "  -0.34%  "

# Row 3: Ethereum
$ws.Range("D3").Value = "1.595.57"
$ws.Range("E3").Value = "  -0.28%  "

# Row 4: TetherUSD
$ws.Range("E4").Value = "  -0.08%  "

# Row 5: BNB
$ws.Range("D5").Value = "'210.73"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.32%  "

# Row 6: XRP
$ws.Range("E6").Value = "  -0.32%  "

# Row 7: USDC
$ws.Range("E7").Value = "  -0.07%  "

# Row 8: Dogecoin
$ws.Range("E8").Value = "  -0.64%  "

# Row 9: Cardano
$ws.Range("D9").Value = "'0.247"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -0.37%  "

# Row 10: Solana
$ws.Range("D10").Value = "'19.60"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +0.30%  "

# Row 11: TRON
$ws.Range("D11").Value = "'0.0845"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +0.46%  "

# Row 12: WrappedliquidstakedEther2.0
$ws.Range("D12").Value = "1.818.51"
$ws.Range("E12").Value = "  -0.36%  "

# Row 13: WrappedEther
$ws.Range("D13").Value = "1.603.77"
$ws.Range("E13").Value = "  +0.58%  "

# Row 14: Polkadot
$ws.Range("E14").Value = "  +0.00%  "

# Row 15: Polygon
$ws.Range("E15").Value = "  -0.26%  "

# Row 16: Litecoin
$ws.Range("D16").Value = "'64.51"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -1.23%  "

# Row 17: WrappedBTC
$ws.Range("D17").Value = "26.588.78"

# Row 18: ShibaInu
$ws.Range("E18").Value = "  -2.35%  "

# Row 19: Dai
$ws.Range("E19").Value = "  -0.06%  "

# Row 20: BitcoinCash
$ws.Range("D20").Value = "'208.97"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -0.06%  "

# Row 21: Chainlink
$ws.Range("E21").Value = "  -1.89%  "

# Row 22: Uniswap
$ws.Range("E22").Value = "  +0.13%  "

# Row 23: Toncoin
$ws.Range("D23").Value = "'2.22"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -3.60%  "

# Row 24: Avalanche
$ws.Range("D24").Value = "'8.96"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +0.15%  "

# Row 25: Monero
$ws.Range("D25").Value = "'145.03"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +1.75%  "

# Row 26: BinanceUSD
$ws.Range("D26").Value = "'1.00"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -0.05%  "

# Row 27: Cosmos
$ws.Range("E27").Value = "  +0.13%  "

# Row 28: Stellar
$ws.Range("E28").Value = "  -0.61%  "

# Row 29: EthereumClassic
$ws.Range("D29").Value = "'15.30"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -0.29%  "

# Row 30: Hedera
$ws.Range("E30").Value = "  -2.60%  "

# Row 31: PancakeSwap
$ws.Range("E31").Value = "  -0.41%  "

# Row 32: Filecoin
$ws.Range("D32").Value = "'3.26"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +0.06%  "

# Row 33: InternetComputer(DFINITY)
$ws.Range("D33").Value = "'2.97"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -0.27%  "

# Row 34: Maker
$ws.Range("D34").Value = "1.282.53"
$ws.Range("E34").Value = "  -0.75%  "

# Row 35: HuobiToken
$ws.Range("E35").Value = "  +0.35%  "

# Row 36: WEMIXToken
$ws.Range("D36").Value = "'1.22"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +11.83%  "

# Row 37: LidoDAOToken -> ImmutableX
$ws.Range("B37").Value = "ImmutableX"
$ws.Range("C37").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D37").Value = "'0.600"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -3.39%  "

# Row 38: ImmutableX -> LidoDAOToken
$ws.Range("B38").Value = "LidoDAOToken"
$ws.Range("C38").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("D38").Value = "'1.48"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -0.96%  "

# Row 39: VeChain
$ws.Range("E39").Value = "  -1.76%  "

# Row 40: ARBITRUM
$ws.Range("E40").Value = "  -0.49%  "

# Row 41: FraxShare
$ws.Range("E41").Value = "  +0.39%  "

# Row 42: MXToken
$ws.Range("D42").Value = "'2.16"
$ws.Range("D42").Style = "Normal"

# Row 43: TrustWalletToken
$ws.Range("D43").Value = "'0.769"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -2.02%  "

# Row 44: Aave
$ws.Range("D44").Value = "'62.75"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -0.80%  "

# Row 45: RocketPoolETH
$ws.Range("D45").Value = "1.730.88"
$ws.Range("E45").Value = "  -0.37%  "

# Row 46: Quant
$ws.Range("D46").Value = "'89.40"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -2.09%  "

# Row 47: RenderToken
$ws.Range("E47").Value = "  +0.06%  "

# Row 48: Algorand
$ws.Range("E48").Value = "  +2.51%  "

# Row 49: Cronos -> BabyDogeCoin
$ws.Range("B49").Value = "BabyDogeCoin"
$ws.Range("C49").Value = "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
$ws.Range("D49").Value = "0.0₆0103"
$ws.Range("E49").Value = "  -2.95%  "

# Row 50: USDD -> Cronos
$ws.Range("B50").Value = "Cronos"
$ws.Range("C50").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D50").Value = "'0.0513"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +0.47%  "

# Row 51: EnergySwap -> USDD
$ws.Range("B51").Value = "USDD"
$ws.Range("C51").Value = "https://coinranking.com/coin/z2PZIKQL7+usdd-usdd"
$ws.Range("D51").Value = "'1.00"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -0.07%  "
